$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1007
$ws.Range("D2").Value = 3512
$ws.Range("E2").Value = 28.67
$ws.Range("F2").Value = 27.18
$ws.Range("G2").Value = 30.17
$ws.Range("C3").Value = 740
$ws.Range("D3").Value = 3512
$ws.Range("E3").Value = 21.07
$ws.Range("F3").Value = 19.72
$ws.Range("G3").Value = 22.42
$ws.Range("B4").Value = 'Symptom – Nervous'
$ws.Range("C4").Value = 367
$ws.Range("D4").Value = 3512
$ws.Range("E4").Value = 10.45
$ws.Range("F4").Value = 9.44
$ws.Range("G4").Value = 11.46
$ws.Range("B5").Value = 'Symptom – Digestive'
$ws.Range("C5").Value = 342
$ws.Range("D5").Value = 3512
$ws.Range("E5").Value = 9.74
$ws.Range("F5").Value = 8.76
$ws.Range("G5").Value = 10.72
$ws.Range("B6").Value = 'Diseases (patient-stated)'
$ws.Range("C6").Value = 240
$ws.Range("D6").Value = 3512
$ws.Range("E6").Value = 6.83
$ws.Range("F6").Value = 6
$ws.Range("G6").Value = 7.67
$ws.Range("B7").Value = 'Symptom – Circulatory'
$ws.Range("C7").Value = 227
$ws.Range("D7").Value = 3512
$ws.Range("E7").Value = 6.46
$ws.Range("F7").Value = 5.65
$ws.Range("G7").Value = 7.28
$ws.Range("B8").Value = 'Other'
$ws.Range("C8").Value = 196
$ws.Range("D8").Value = 3512
$ws.Range("E8").Value = 5.58
$ws.Range("F8").Value = 4.82
$ws.Range("G8").Value = 6.34
$ws.Range("B9").Value = 'Uncodable/Unknown'
$ws.Range("C9").Value = 138
$ws.Range("D9").Value = 3512
$ws.Range("E9").Value = 3.93
$ws.Range("F9").Value = 3.29
$ws.Range("G9").Value = 4.57
$ws.Range("B10").Value = 'Symptom – General'
$ws.Range("C10").Value = 97
$ws.Range("D10").Value = 3512
$ws.Range("E10").Value = 2.76
$ws.Range("F10").Value = 2.22
$ws.Range("G10").Value = 3.3
$ws.Range("B11").Value = 'Administrative'
$ws.Range("C11").Value = 80
$ws.Range("D11").Value = 3512
$ws.Range("E11").Value = 2.28
$ws.Range("F11").Value = 1.78
$ws.Range("G11").Value = 2.77
$ws.Range("B12").Value = 'Symptom – Skin/Hair/Nails'
$ws.Range("C12").Value = 78
$ws.Range("D12").Value = 3512
$ws.Range("E12").Value = 2.22
$ws.Range("F12").Value = 1.73
$ws.Range("G12").Value = 2.71
$ws.Range("C13").Value = 989
$ws.Range("E13").Value = 49.87
$ws.Range("F13").Value = 47.67
$ws.Range("G13").Value = 52.07
$ws.Range("C18").Value = 114
$ws.Range("E18").Value = 5.75
$ws.Range("F18").Value = 4.72
$ws.Range("G18").Value = 6.77
$ws.Range("C20").Value = 59
$ws.Range("E20").Value = 2.98
$ws.Range("F20").Value = 2.23
$ws.Range("G20").Value = 3.72
$ws.Range("B21").Value = 'Administrative'
$ws.Range("C21").Value = 30
$ws.Range("E21").Value = 1.51
$ws.Range("F21").Value = 0.98
$ws.Range("G21").Value = 2.05
$ws.Range("B22").Value = 'Symptom – Skin/Hair/Nails'
$ws.Range("C22").Value = 27
$ws.Range("E22").Value = 1.36
$ws.Range("F22").Value = 0.85
$ws.Range("G22").Value = 1.87
$ws.Range("C23").Value = 25
$ws.Range("E23").Value = 1.26
$ws.Range("F23").Value = 0.77
$ws.Range("G23").Value = 1.75
$ws.Range("A24").Value = 'pco2_threshold_any'
$ws.Range("C24").Value = 2037
$ws.Range("D24").Value = 6246
$ws.Range("E24").Value = 32.61
$ws.Range("F24").Value = 31.45
$ws.Range("G24").Value = 33.78
$ws.Range("A25").Value = 'pco2_threshold_any'
$ws.Range("B25").Value = 'Injuries & adverse effects'
$ws.Range("C25").Value = 1030
$ws.Range("D25").Value = 6246
$ws.Range("E25").Value = 16.49
$ws.Range("F25").Value = 15.57
$ws.Range("G25").Value = 17.41
$ws.Range("A26").Value = 'pco2_threshold_any'
$ws.Range("C26").Value = 739
$ws.Range("D26").Value = 6246
$ws.Range("E26").Value = 11.83
$ws.Range("F26").Value = 11.03
$ws.Range("G26").Value = 12.63
$ws.Range("A27").Value = 'pco2_threshold_any'
$ws.Range("B27").Value = 'Symptom – Digestive'
$ws.Range("C27").Value = 617
$ws.Range("D27").Value = 6246
$ws.Range("E27").Value = 9.880000000000001
$ws.Range("F27").Value = 9.140000000000001
$ws.Range("G27").Value = 10.62
$ws.Range("A28").Value = 'pco2_threshold_any'
$ws.Range("C28").Value = 473
$ws.Range("D28").Value = 6246
$ws.Range("E28").Value = 7.57
$ws.Range("F28").Value = 6.92
$ws.Range("G28").Value = 8.23
$ws.Range("A29").Value = 'pco2_threshold_any'
$ws.Range("C29").Value = 378
$ws.Range("D29").Value = 6246
$ws.Range("E29").Value = 6.05
$ws.Range("F29").Value = 5.46
$ws.Range("G29").Value = 6.64
$ws.Range("A30").Value = 'pco2_threshold_any'
$ws.Range("B30").Value = 'Diseases (patient-stated)'
$ws.Range("C30").Value = 358
$ws.Range("D30").Value = 6246
$ws.Range("E30").Value = 5.73
$ws.Range("F30").Value = 5.16
$ws.Range("G30").Value = 6.31
$ws.Range("A31").Value = 'pco2_threshold_any'
$ws.Range("B31").Value = 'Symptom – General'
$ws.Range("C31").Value = 212
$ws.Range("D31").Value = 6246
$ws.Range("E31").Value = 3.39
$ws.Range("F31").Value = 2.95
$ws.Range("G31").Value = 3.84
$ws.Range("A32").Value = 'pco2_threshold_any'
$ws.Range("B32").Value = 'Uncodable/Unknown'
$ws.Range("C32").Value = 174
$ws.Range("D32").Value = 6246
$ws.Range("E32").Value = 2.79
$ws.Range("F32").Value = 2.38
$ws.Range("G32").Value = 3.19
$ws.Range("A33").Value = 'pco2_threshold_any'
$ws.Range("B33").Value = 'Administrative'
$ws.Range("C33").Value = 121
$ws.Range("D33").Value = 6246
$ws.Range("E33").Value = 1.94
$ws.Range("F33").Value = 1.6
$ws.Range("G33").Value = 2.28
$ws.Range("A34").Value = 'pco2_threshold_any'
$ws.Range("B34").Value = 'Symptom – Skin/Hair/Nails'
$ws.Range("C34").Value = 107
$ws.Range("D34").Value = 6246
$ws.Range("E34").Value = 1.71
$ws.Range("F34").Value = 1.39
$ws.Range("G34").Value = 2.03
$ws.Range("A35").Value = 'unknown_hypercap_threshold'
$ws.Range("C35").Value = 117
$ws.Range("D35").Value = 396
$ws.Range("E35").Value = 29.55
$ws.Range("F35").Value = 25.05
$ws.Range("G35").Value = 34.04
$ws.Range("A36").Value = 'unknown_hypercap_threshold'
$ws.Range("C36").Value = 65
$ws.Range("D36").Value = 396
$ws.Range("E36").Value = 16.41
$ws.Range("F36").Value = 12.77
$ws.Range("G36").Value = 20.06
$ws.Range("A37").Value = 'unknown_hypercap_threshold'
$ws.Range("B37").Value = 'Symptom – Nervous'
$ws.Range("C37").Value = 46
$ws.Range("D37").Value = 396
$ws.Range("E37").Value = 11.62
$ws.Range("F37").Value = 8.460000000000001
$ws.Range("G37").Value = 14.77
$ws.Range("A38").Value = 'unknown_hypercap_threshold'
$ws.Range("B38").Value = 'Symptom – Digestive'
$ws.Range("C38").Value = 40
$ws.Range("D38").Value = 396
$ws.Range("E38").Value = 10.1
$ws.Range("F38").Value = 7.13
$ws.Range("G38").Value = 13.07
$ws.Range("A39").Value = 'unknown_hypercap_threshold'
$ws.Range("C39").Value = 36
$ws.Range("D39").Value = 396
$ws.Range("E39").Value = 9.09
$ws.Range("F39").Value = 6.26
$ws.Range("G39").Value = 11.92
$ws.Range("A40").Value = 'unknown_hypercap_threshold'
$ws.Range("C40").Value = 28
$ws.Range("D40").Value = 396
$ws.Range("E40").Value = 7.07
$ws.Range("F40").Value = 4.55
$ws.Range("G40").Value = 9.6
$ws.Range("A41").Value = 'unknown_hypercap_threshold'
$ws.Range("C41").Value = 22
$ws.Range("D41").Value = 396
$ws.Range("E41").Value = 5.56
$ws.Range("F41").Value = 3.3
$ws.Range("G41").Value = 7.81
$ws.Range("A42").Value = 'unknown_hypercap_threshold'
$ws.Range("C42").Value = 15
$ws.Range("D42").Value = 396
$ws.Range("E42").Value = 3.79
$ws.Range("F42").Value = 1.91
$ws.Range("G42").Value = 5.67
$ws.Range("A43").Value = 'unknown_hypercap_threshold'
$ws.Range("B43").Value = 'Administrative'
$ws.Range("C43").Value = 9
$ws.Range("D43").Value = 396
$ws.Range("E43").Value = 2.27
$ws.Range("F43").Value = 0.8
$ws.Range("G43").Value = 3.74
$ws.Range("A44").Value = 'unknown_hypercap_threshold'
$ws.Range("B44").Value = 'Symptom – Skin/Hair/Nails'
$ws.Range("C44").Value = 9
$ws.Range("D44").Value = 396
$ws.Range("E44").Value = 2.27
$ws.Range("F44").Value = 0.8
$ws.Range("G44").Value = 3.74
$ws.Range("A45").Value = 'unknown_hypercap_threshold'
$ws.Range("C45").Value = 9
$ws.Range("D45").Value = 396
$ws.Range("E45").Value = 2.27
$ws.Range("F45").Value = 0.8
$ws.Range("G45").Value = 3.74
$ws.Range("C46").Value = 1380
$ws.Range("D46").Value = 3507
$ws.Range("E46").Value = 39.35
$ws.Range("F46").Value = 37.73
$ws.Range("G46").Value = 40.97
$ws.Range("C47").Value = 460
$ws.Range("D47").Value = 3507
$ws.Range("E47").Value = 13.12
$ws.Range("F47").Value = 12
$ws.Range("G47").Value = 14.23
$ws.Range("B48").Value = 'Injuries & adverse effects'
$ws.Range("C48").Value = 419
$ws.Range("D48").Value = 3507
$ws.Range("E48").Value = 11.95
$ws.Range("F48").Value = 10.87
$ws.Range("G48").Value = 13.02
$ws.Range("B49").Value = 'Symptom – Digestive'
$ws.Range("C49").Value = 323
$ws.Range("D49").Value = 3507
$ws.Range("E49").Value = 9.210000000000001
$ws.Range("F49").Value = 8.25
$ws.Range("G49").Value = 10.17
$ws.Range("C50").Value = 291
$ws.Range("D50").Value = 3507
$ws.Range("E50").Value = 8.300000000000001
$ws.Range("F50").Value = 7.38
$ws.Range("G50").Value = 9.210000000000001
$ws.Range("C51").Value = 205
$ws.Range("D51").Value = 3507
$ws.Range("E51").Value = 5.85
$ws.Range("F51").Value = 5.07
$ws.Range("G51").Value = 6.62
$ws.Range("B52").Value = 'Diseases (patient-stated)'
$ws.Range("C52").Value = 156
$ws.Range("D52").Value = 3507
$ws.Range("E52").Value = 4.45
$ws.Range("F52").Value = 3.77
$ws.Range("G52").Value = 5.13
$ws.Range("B53").Value = 'Symptom – General'
$ws.Range("C53").Value = 129
$ws.Range("D53").Value = 3507
$ws.Range("E53").Value = 3.68
$ws.Range("F53").Value = 3.06
$ws.Range("G53").Value = 4.3
$ws.Range("B54").Value = 'Uncodable/Unknown'
$ws.Range("C54").Value = 57
$ws.Range("D54").Value = 3507
$ws.Range("E54").Value = 1.63
$ws.Range("F54").Value = 1.21
$ws.Range("G54").Value = 2.04
$ws.Range("B55").Value = 'Administrative'
$ws.Range("C55").Value = 54
$ws.Range("D55").Value = 3507
$ws.Range("E55").Value = 1.54
$ws.Range("F55").Value = 1.13
$ws.Range("G55").Value = 1.95
$ws.Range("B56").Value = 'Symptom – Skin/Hair/Nails'
$ws.Range("C56").Value = 33
$ws.Range("D56").Value = 3507
$ws.Range("E56").Value = 0.9399999999999999
$ws.Range("F56").Value = 0.62
$ws.Range("G56").Value = 1.26
